$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price column keeps exact text formatting (avoid Excel auto-converting numeric-looking strings)
$ws.Range("D2:D17").NumberFormat = "@"
$ws.Range("D19:D26").NumberFormat = "@"
$ws.Range("D28:D40").NumberFormat = "@"
$ws.Range("D42:D46").NumberFormat = "@"
$ws.Range("D48:D51").NumberFormat = "@"

$ws.Range("D2").Value = "65.934.07"
$ws.Range("E2").Value = "  +0.36%  "

$ws.Range("D3").Value = "3.317.17"
$ws.Range("E3").Value = "  +1.55%  "

$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "560.59"
$ws.Range("E5").Value = "  +0.82%  "

$ws.Range("D6").Value = "185.18"
$ws.Range("E6").Value = "  +0.73%  "

$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.12%  "

$ws.Range("D8").Value = "3.309.68"
$ws.Range("E8").Value = "  +1.55%  "

$ws.Range("D9").Value = "0.572"
$ws.Range("E9").Value = "  -2.71%  "

$ws.Range("D10").Value = "0.175"
$ws.Range("E10").Value = "  -4.46%  "

$ws.Range("D11").Value = "0.576"
$ws.Range("E11").Value = "  -0.99%  "

$ws.Range("D12").Value = "45.84"
$ws.Range("E12").Value = "  -2.11%  "

$ws.Range("D13").Value = "0.0000263"
$ws.Range("E13").Value = "  -0.55%  "

$ws.Range("D14").Value = "3.842.01"
$ws.Range("E14").Value = "  +1.45%  "

$ws.Range("D15").Value = "8.43"
$ws.Range("E15").Value = "  -1.92%  "

$ws.Range("D16").Value = "580.71"
$ws.Range("E16").Value = "  -8.43%  "

$ws.Range("D17").Value = "65.929.55"
$ws.Range("E17").Value = "  +0.50%  "

$ws.Range("D19").Value = "3.305.24"
$ws.Range("E19").Value = "  +1.39%  "

$ws.Range("D20").Value = "17.58"
$ws.Range("E20").Value = "  -2.24%  "

$ws.Range("D21").Value = "10.88"
$ws.Range("E21").Value = "  -3.67%  "

$ws.Range("D22").Value = "0.891"
$ws.Range("E22").Value = "  -0.99%  "

$ws.Range("D23").Value = "17.68"
$ws.Range("E23").Value = "  -2.94%  "

$ws.Range("D24").Value = "4.99"
$ws.Range("E24").Value = "  +2.24%  "

$ws.Range("D25").Value = "97.90"
$ws.Range("E25").Value = "  -8.39%  "

$ws.Range("D26").Value = "3.94"
$ws.Range("E26").Value = "  -0.11%  "

$ws.Range("E27").Value = "  -0.69%  "

$ws.Range("D28").Value = "2.69"
$ws.Range("E28").Value = "  +1.07%  "

$ws.Range("D29").Value = "9.33"
$ws.Range("E29").Value = "  -1.94%  "

$ws.Range("D30").Value = "8.44"
$ws.Range("E30").Value = "  -2.15%  "

$ws.Range("D31").Value = "30.58"
$ws.Range("E31").Value = "  +1.49%  "

$ws.Range("D32").Value = "6.64"
$ws.Range("E32").Value = "  +6.60%  "

$ws.Range("D33").Value = "3.70"
$ws.Range("E33").Value = "  -5.90%  "

$ws.Range("D34").Value = "558.78"
$ws.Range("E34").Value = "  +7.09%  "

$ws.Range("D35").Value = "10.79"
$ws.Range("E35").Value = "  -1.90%  "

$ws.Range("B36").Value = "Maker"
$ws.Range("C36").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D36").Value = "3.747.12"
$ws.Range("E36").Value = "  +0.15%  "

$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").Value = "0.103"
$ws.Range("E37").Value = "  -1.46%  "

$ws.Range("D38").Value = "1.00"
$ws.Range("E38").Value = "  +0.14%  "

$ws.Range("D39").Value = "55.62"
$ws.Range("E39").Value = "  -3.33%  "

$ws.Range("D40").Value = "33.27"
$ws.Range("E40").Value = "  +1.54%  "

$ws.Range("E41").Value = "  -3.26%  "

$ws.Range("D42").Value = "3.14"
$ws.Range("E42").Value = "  -6.59%  "

$ws.Range("D43").Value = "0.0₃0682"
$ws.Range("E43").Value = "  -7.07%  "

$ws.Range("B44").Value = "Fetch.AI"
$ws.Range("C44").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D44").Value = "2.58"
$ws.Range("E44").Value = "  -4.33%  "

$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").Value = "3.32"
$ws.Range("E45").Value = "  +4.30%  "

$ws.Range("D46").Value = "0.332"
$ws.Range("E46").Value = "  -0.80%  "

$ws.Range("E47").Value = "  -8.23%  "

$ws.Range("D48").Value = "0.0408"
$ws.Range("E48").Value = "  -0.93%  "

$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").Value = "0.126"
$ws.Range("E49").Value = "  -1.98%  "

$ws.Range("B50").Value = "FirstDigitalUSD"
$ws.Range("C50").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D50").Value = "0.999"
$ws.Range("E50").Value = "  +0.11%  "

$ws.Range("D51").Value = "2.50"
$ws.Range("E51").Value = "  -3.47%  "
